$d = $word.ActiveDocument

function Get-ParaIndexContaining($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Drop the trailing "Тайна Христова, Творение, ... Тьма" index-line
#    paragraph (whole paragraph removed).
# ------------------------------------------------------------------
$idx = Get-ParaIndexContaining $d "Тайна Христова, Творение"
if ($idx -gt 0) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# ------------------------------------------------------------------
# 2) Drop the "This PDF version is provided under the same license."
#    paragraph entirely.
# ------------------------------------------------------------------
$idx = Get-ParaIndexContaining $d "This PDF version is provided under the same license."
if ($idx -gt 0) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# ------------------------------------------------------------------
# 3) Rewrite the license paragraph (bold title + description) into the
#    new Biblica Study Notes copy, preserving the leading empty run and
#    the bold title run's formatting.
# ------------------------------------------------------------------
$idx = Get-ParaIndexContaining $d "Ключевые термины (Biblica)"
$idx2 = Get-ParaIndexContaining $d " (Russian) is based on"
if ($idx2 -gt 0) { $idx = $idx2 }
$para = $d.Paragraphs.Item($idx)
$rng = $para.Range

$newParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr/><w:r><w:rPr><w:lang w:val="ru_RU" w:bidi="ru_RU"/></w:rPr></w:r><w:r><w:rPr><w:b/><w:lang w:val="ru_RU" w:bidi="ru_RU"/></w:rPr><w:t>Biblica Study Notes (Key Terms)</w:t></w:r><w:r><w:rPr><w:lang w:val="ru_RU" w:bidi="ru_RU"/></w:rPr><w:t xml:space="preserve"> &#169; 2023 Biblica Inc. Released under CC BY-SA 4.0 license. </w:t></w:r><w:r><w:rPr><w:lang w:val="ru_RU" w:bidi="ru_RU"/></w:rPr><w:t>Biblica Study Notes</w:t></w:r><w:r><w:rPr><w:lang w:val="ru_RU" w:bidi="ru_RU"/></w:rPr><w:t xml:space="preserve"> has been adapted in the following languages: Tok Pisin, Arabic (&#1593;&#1585;&#1576;&#1610;), French (Fran&#231;ais), Hindi (&#2361;&#2367;&#2306;&#2342;&#2368;), Indonesian (Bahasa Indonesia), Portuguese (Portugu&#234;s), Russian (&#1056;&#1091;&#1089;&#1089;&#1082;&#1080;&#1081;), Spanish (Espa&#241;ol), Swahili (Kiswahili), and Simplified Chinese (&#31616;&#20307;&#20013;&#25991;)from Biblica Study Notes &#169; 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual.</w:t></w:r><w:r><w:rPr><w:lang w:val="ru_RU" w:bidi="ru_RU"/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($newParaXml)

# ------------------------------------------------------------------
# 4) Drop the "License Information" heading paragraph entirely.
# ------------------------------------------------------------------
$idx = Get-ParaIndexContaining $d "License Information"
if ($idx -gt 0) {
    $d.Paragraphs.Item($idx).Range.Delete()
}
